$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.946.10"
$ws.Cells.Item(2, 5).Value = "  +4.54%  "

$ws.Cells.Item(3, 4).Value = "1.781.16"
$ws.Cells.Item(3, 5).Value = "  +3.13%  "

$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.36%  "

$ws.Cells.Item(5, 4).Value = "'244.03"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.85%  "

$ws.Cells.Item(6, 4).Value = "'1.0000"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.28%  "

$ws.Cells.Item(7, 4).Value = "'0.4911"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.18%  "

$ws.Cells.Item(8, 4).Value = "'0.2671"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.11%  "

$ws.Cells.Item(9, 4).Value = "'0.06254"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.40%  "

$ws.Cells.Item(10, 4).Value = "1.785.91"
$ws.Cells.Item(10, 5).Value = "  +3.41%  "

$ws.Cells.Item(11, 4).Value = "'16.33"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.00%  "

$ws.Cells.Item(12, 4).Value = "'0.07044"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.65%  "

$ws.Cells.Item(13, 4).Value = "'0.6262"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.60%  "

$ws.Cells.Item(14, 4).Value = "'4.629"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.79%  "

$ws.Cells.Item(15, 4).Value = "'80.01"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.62%  "

$ws.Cells.Item(16, 5).Value = "  +0.35%  "

$ws.Cells.Item(17, 4).Value = "27.922.83"
$ws.Cells.Item(17, 5).Value = "  +5.29%  "

$ws.Cells.Item(18, 4).Value = "'0.9998"
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'0.000007212"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.08%  "

$ws.Cells.Item(20, 4).Value = "'11.93"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +4.51%  "

$ws.Cells.Item(21, 4).Value = "2.007.79"
$ws.Cells.Item(21, 5).Value = "  +2.98%  "

$ws.Cells.Item(22, 4).Value = "'4.578"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +3.05%  "

$ws.Cells.Item(23, 5).Value = "  +1.34%  "

$ws.Cells.Item(24, 4).Value = "'5.231"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.41%  "

$ws.Cells.Item(25, 4).Value = "'141.53"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.47%  "

$ws.Cells.Item(26, 4).Value = "'15.70"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.10%  "

$ws.Cells.Item(27, 4).Value = "'1.860"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +6.46%  "

$ws.Cells.Item(28, 4).Value = "'109.27"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.71%  "

$ws.Cells.Item(29, 4).Value = "'1.397"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.41%  "

$ws.Cells.Item(30, 4).Value = "'4.205"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +7.30%  "

$ws.Cells.Item(31, 4).Value = "'0.08277"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.59%  "

$ws.Cells.Item(32, 4).Value = "'3.793"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.42%  "

$ws.Cells.Item(33, 4).Value = "'0.04859"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +7.98%  "

$ws.Cells.Item(34, 5).Value = "  +6.62%  "

$ws.Cells.Item(35, 5).Value = "  +0.18%  "

$ws.Cells.Item(36, 4).Value = "'0.6494"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +3.65%  "

$ws.Cells.Item(37, 4).Value = "'0.9466"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.01%  "

$ws.Cells.Item(38, 4).Value = "'2.583"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +6.60%  "

$ws.Cells.Item(39, 5).Value = "  +1.67%  "

$ws.Cells.Item(40, 4).Value = "'5.947"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +6.58%  "

$ws.Cells.Item(41, 4).Value = "'0.01551"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.45%  "

$ws.Cells.Item(42, 4).Value = "'1.000"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.32%  "

$ws.Cells.Item(43, 4).Value = "'100.00"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.43%  "

$ws.Cells.Item(44, 4).Value = "'0.3980"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.08%  "

$ws.Cells.Item(45, 4).Value = "'7.165"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +3.75%  "

$ws.Cells.Item(46, 4).Value = "'0.1200"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +3.52%  "

$ws.Cells.Item(47, 4).Value = "'0.05417"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.64%  "

$ws.Cells.Item(48, 4).Value = "'7.973"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.92%  "

$ws.Cells.Item(49, 4).Value = "'1.295"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +4.78%  "

$ws.Cells.Item(50, 4).Value = "'30.68"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.05%  "

$ws.Cells.Item(51, 4).Value = "'52.93"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.35%  "
